$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 452 ("「どうやらこの好奇心旺盛なペンギンは…」" entry) entirely.
# This shifts all subsequent rows up by one, matching the author's edit
# (which removed that post and renumbered everything below it).
$ws.Rows.Item(452).Delete()
